$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B2").Value = 0.6304347826086957
$ws.Range("C2").Value = 0.7631578947368421
$ws.Range("D2").Value = 0.6904761904761905
$ws.Range("B3").Value = 0.8085106382978723
$ws.Range("C3").Value = 0.6909090909090909
$ws.Range("D3").Value = 0.7450980392156863
$ws.Range("B4").Value = 0.7204301075268817
$ws.Range("C4").Value = 0.7204301075268817
$ws.Range("D4").Value = 0.7204301075268817
$ws.Range("E4").Value = 0.7204301075268817
$ws.Range("B5").Value = 0.719472710453284
$ws.Range("C5").Value = 0.7270334928229665
$ws.Range("D5").Value = 0.7177871148459384
$ws.Range("B6").Value = 0.7357484607044452
$ws.Range("C6").Value = 0.7204301075268817
$ws.Range("D6").Value = 0.722779434354387
$ws.Range("B7").Value = 0.6
$ws.Range("C7").Value = 0.7105263157894737
$ws.Range("D7").Value = 0.6506024096385543
$ws.Range("B8").Value = 0.7708333333333334
$ws.Range("C8").Value = 0.6727272727272727
$ws.Range("D8").Value = 0.7184466019417476
$ws.Range("B9").Value = 0.6881720430107527
$ws.Range("C9").Value = 0.6881720430107527
$ws.Range("D9").Value = 0.6881720430107527
$ws.Range("E9").Value = 0.6881720430107527
$ws.Range("B10").Value = 0.6854166666666667
$ws.Range("C10").Value = 0.6916267942583731
$ws.Range("D10").Value = 0.684524505790151
$ws.Range("B11").Value = 0.7010304659498209
$ws.Range("C11").Value = 0.6881720430107527
$ws.Range("D11").Value = 0.690725319065174
$ws.Range("B17").Value = 0.7777777777777778
$ws.Range("C17").Value = 0.7368421052631579
$ws.Range("D17").Value = 0.7567567567567567
$ws.Range("B18").Value = 0.8245614035087719
$ws.Range("C18").Value = 0.8545454545454545
$ws.Range("D18").Value = 0.8392857142857144
$ws.Range("B19").Value = 0.8064516129032258
$ws.Range("C19").Value = 0.8064516129032258
$ws.Range("D19").Value = 0.8064516129032258
$ws.Range("E19").Value = 0.8064516129032258
$ws.Range("B20").Value = 0.8011695906432749
$ws.Range("C20").Value = 0.7956937799043062
$ws.Range("D20").Value = 0.7980212355212355
$ws.Range("B21").Value = 0.8054455134251399
$ws.Range("C21").Value = 0.8064516129032258
$ws.Range("D21").Value = 0.8055642047577533
$ws.Range("B22").Value = 0.5833333333333334
$ws.Range("C22").Value = 0.7368421052631579
$ws.Range("D22").Value = 0.6511627906976745
$ws.Range("B23").Value = 0.7777777777777778
$ws.Range("C23").Value = 0.6363636363636364
$ws.Range("D23").Value = 0.7000000000000001
$ws.Range("B24").Value = 0.6774193548387096
$ws.Range("C24").Value = 0.6774193548387096
$ws.Range("D24").Value = 0.6774193548387096
$ws.Range("E24").Value = 0.6774193548387096
$ws.Range("B25").Value = 0.6805555555555556
$ws.Range("C25").Value = 0.6866028708133971
$ws.Range("D25").Value = 0.6755813953488372
$ws.Range("B26").Value = 0.6983273596176822
$ws.Range("C26").Value = 0.6774193548387096
$ws.Range("D26").Value = 0.6800450112528132
